$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 457
$ws1.Range("F5").Value = 333
$ws1.Range("F6").Value = 475
$ws1.Range("F8").Value = 2207
$ws1.Range("F9").Value = 52
$ws1.Range("F10").Value = 62
$ws1.Range("F11").Value = 1651
$ws1.Range("F12").Value = 1651
$ws1.Range("F13").Value = 1367
$ws1.Range("F15").Value = 1420
$ws1.Range("F17").Value = 20
$ws1.Range("F18").Value = 589
$ws1.Range("F19").Value = 163
$ws1.Range("F20").Value = 18
$ws1.Range("F21").Value = 7291
$ws1.Range("F22").Value = 8069
$ws1.Range("F25").Value = 203
$ws1.Range("F31").Value = 256
$ws1.Range("F35").Value = 1450
$ws1.Range("F36").Value = 211
$ws1.Range("F37").Value = 229
$ws1.Range("F38").Value = 14
$ws1.Range("F39").Value = 294
$ws1.Range("F40").Value = 14
$ws1.Range("F41").Value = 741
$ws1.Range("F44").Value = 347
$ws1.Range("F48").Value = 181
$ws1.Range("F49").Value = 165

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 457
$ws4.Range("F8").Value = 333
$ws4.Range("F10").Value = 475
$ws4.Range("F11").Value = 52
$ws4.Range("F12").Value = 62
$ws4.Range("F13").Value = 1651
$ws4.Range("F14").Value = 1651
$ws4.Range("F16").Value = 1367
$ws4.Range("F19").Value = 20
$ws4.Range("F20").Value = 589
$ws4.Range("F21").Value = 163
$ws4.Range("F23").Value = 18
$ws4.Range("F24").Value = 7291
$ws4.Range("F25").Value = 8071
$ws4.Range("F27").Value = 203
$ws4.Range("F32").Value = 1450
$ws4.Range("F33").Value = 211
$ws4.Range("F34").Value = 229
$ws4.Range("F35").Value = 14
$ws4.Range("F37").Value = 294
$ws4.Range("F39").Value = 741
$ws4.Range("F44").Value = 347
$ws4.Range("F47").Value = 181
$ws4.Range("F48").Value = 165
